$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Correct the tiny floating point value in A15 (same timestamp, adjusted precision)
$ws.Cells.Item(15, 1).Value = 45876.54187708333

# Add new row 16 with new sensor reading
$ws.Cells.Item(16, 1).Value = 45876.58349000412
$ws.Cells.Item(16, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(16, 2).Value = 2025
$ws.Cells.Item(16, 3).Value = 28
$ws.Cells.Item(16, 4).Value = 20.89
$ws.Cells.Item(16, 5).Value = 74.63
$ws.Cells.Item(16, 6).Value = 68.55
$ws.Cells.Item(16, 7).Value = 13.64
$ws.Cells.Item(16, 8).Value = "ESE"
$ws.Cells.Item(16, 9).Value = 0
$ws.Cells.Item(16, 10).Value = "14:00:13"
